$d = $word.ActiveDocument

$replacements = @(
    @("42×60=", "34×85="),
    @("89×82=", "84×31="),
    @("93×78=", "91×30="),
    @("19×89=", "98×79="),
    @("48×20=", "59×80="),
    @("35×86=", "34×40="),
    @("68×18=", "92×16="),
    @("49×65=", "88×95="),
    @("88×63=", "60×83="),
    @("97×65=", "67×84="),
    @("33×30=", "64×52="),
    @("98×81=", "53×29="),
    @("11×97=", "38×51="),
    @("98×90=", "26×94="),
    @("92×94=", "59×96="),
    @("28×49=", "42×76="),
    @("54×61=", "45×25="),
    @("17×38=", "38×93="),
    @("75×58=", "62×85="),
    @("24×53=", "47×53="),
    @("34×77=", "51×67="),
    @("24×42=", "61×49="),
    @("37×33=", "29×45="),
    @("14×75=", "19×46="),
    @("50×20=", "59×87=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
